{"js": "// Lattice-multiplication worksheet refresh: each table cell holds one\n// exercise rendered as 5 lines (separated by <w:br/>, i.e. \"\\u000b\" in the\n// Office.js text projection):\n//   line 1: \"AB x CD\"          (the two factors)\n//   line 2: \"  C    D\"         (top-of-lattice digits, from the 2nd factor)\n//   line 3: \"  ----\"           (divider, unchanged)\n//   line 4: \"A|    |\"          (left-of-lattice digit, 1st digit of 1st factor)\n//   line 5: \"B|    |\"          (left-of-lattice digit, 2nd digit of 1st factor)\n//\n// The commit swaps in a new set of factor pairs (same table shape / same\n// number of exercises) and regenerates lines 2, 4 and 5 to match. This maps\n// each cell's CURRENT \"AB x CD\" heading to its replacement heading, then\n// rebuilds the other four lines from that replacement using the same\n// left/top digit layout rule the worksheet generator uses.\n\nconst BR = \"\\u000b\";\n\n// Old heading -> new heading, in document (row-major) order, taken from the\n// authoritative edit.\nconst REPLACEMENTS = {\n  \"55 x 90\": \"76 x 95\",\n  \"98 x 27\": \"31 x 95\",\n  \"83 x 13\": \"23 x 34\",\n  \"52 x 22\": \"54 x 55\",\n  \"78 x 99\": \"29 x 25\",\n  \"15 x 67\": \"41 x 95\",\n  \"74 x 44\": \"89 x 65\",\n  \"47 x 30\": \"27 x 23\",\n  \"22 x 63\": \"62 x 18\",\n  \"42 x 70\": \"92 x 96\",\n  \"35 x 93\": \"69 x 50\",\n  \"99 x 56\": \"98 x 94\",\n  \"63 x 54\": \"27 x 31\",\n  \"30 x 15\": \"77 x 27\",\n  \"10 x 23\": \"93 x 76\"\n};\n\nfunction buildCellText(heading) {\n  const [a, b] = heading.split(\" x \");\n  const top = \"  \" + b[0] + \"    \" + b[1];\n  const left1 = a[0] + \"|    |\";\n  const left2 = a[1] + \"|    |\";\n  return [heading, top, \"  ----\", left1, left2].join(BR);\n}\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst cellParagraphs = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < 3; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    cellParagraphs.push(para);\n  }\n}\nawait context.sync();\n\nfor (const para of cellParagraphs) {\n  const fullText = para.text || \"\";\n  const heading = fullText.split(BR)[0];\n  const newHeading = REPLACEMENTS[heading];\n  if (newHeading) {\n    para.insertText(buildCellText(newHeading), \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet refresh: each table cell holds one\n# exercise rendered as 5 lines (separated by <w:br/>, which COM's\n# Range.Text surfaces as chr(11)/vertical-tab):\n#   line 1: \"AB x CD\"          (the two factors)\n#   line 2: \"  C    D\"         (top-of-lattice digits, from the 2nd factor)\n#   line 3: \"  ----\"           (divider, unchanged)\n#   line 4: \"A|    |\"          (left-of-lattice digit, 1st digit of 1st factor)\n#   line 5: \"B|    |\"          (left-of-lattice digit, 2nd digit of 1st factor)\n#\n# The commit swaps in a new set of factor pairs (same table shape / same\n# number of exercises) and regenerates lines 2, 4 and 5 to match. This maps\n# each cell's CURRENT \"AB x CD\" heading to its replacement heading, then\n# rebuilds the other four lines from that replacement using the same\n# left/top digit layout rule the worksheet generator uses.\n\n$d = $word.ActiveDocument\n\n$BR = [char]11\n\n$replacements = @{\n    \"55 x 90\" = \"76 x 95\"\n    \"98 x 27\" = \"31 x 95\"\n    \"83 x 13\" = \"23 x 34\"\n    \"52 x 22\" = \"54 x 55\"\n    \"78 x 99\" = \"29 x 25\"\n    \"15 x 67\" = \"41 x 95\"\n    \"74 x 44\" = \"89 x 65\"\n    \"47 x 30\" = \"27 x 23\"\n    \"22 x 63\" = \"62 x 18\"\n    \"42 x 70\" = \"92 x 96\"\n    \"35 x 93\" = \"69 x 50\"\n    \"99 x 56\" = \"98 x 94\"\n    \"63 x 54\" = \"27 x 31\"\n    \"30 x 15\" = \"77 x 27\"\n    \"10 x 23\" = \"93 x 76\"\n}\n\nfunction Build-CellText($heading) {\n    # NOTE: avoid the `+` operator on digit-only substrings here \u2014 this\n    # host's PowerShell numeric-coercion treats e.g. \"  9    \" + \"5\" as an\n    # ADDITION (-> 14) rather than string concatenation once both sides look\n    # numeric after trimming. `-f` (format) and string interpolation both\n    # stay string-typed, so use those instead.\n    $parts = $heading -split \" x \"\n    $a = $parts[0]\n    $b = $parts[1]\n    $top = \"  {0}    {1}\" -f $b.Substring(0,1), $b.Substring(1,1)\n    $left1 = \"{0}|    |\" -f $a.Substring(0,1)\n    $left2 = \"{0}|    |\" -f $a.Substring(1,1)\n    return ($heading, $top, \"  ----\", $left1, $left2) -join $BR\n}\n\n$table = $d.Tables.Item(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $table.Cell($r, $c)\n        $raw = $cell.Range.Text\n        # Strip Word's end-of-cell marker (CR + BEL, chr(13)+chr(7)).\n        $raw = $raw.TrimEnd([char]7).TrimEnd([char]13)\n        $heading = $raw.Split($BR)[0]\n        if ($replacements.ContainsKey($heading)) {\n            $newHeading = $replacements[$heading]\n            $cell.Range.Text = Build-CellText $newHeading\n        }\n    }\n}\n"}
